$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "606.65" are not
# auto-converted to numbers by Excel (matches original inlineStr encoding).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '66.564.86'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.560.97'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '606.65'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '145.04'
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.493'
$ws.Range("E8").Value = '  +2.68%  '
$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").Value = '8.03'
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.136'
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '0.416'
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '4.162.87'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000207'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '30.03'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.530.50'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '66.573.34'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.115'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '11.50'
$ws.Range("E18").Value = '  +5.16%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '6.22'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '15.01'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '431.87'
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '0.617'
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '79.27'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.700.89'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  -1.85%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '8.03'
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.51'
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '9.17'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("B31").Value = 'RenzoRestakedETH'
$ws.Range("C31").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D31").Value = '3.554.70'
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '1.45'
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '25.34'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.155'
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '7.84'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.72'
$ws.Range("E37").Value = '  -1.66%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '5.62'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '173.32'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.0851'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '5.20'
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '0.889'
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.93'
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.51'
$ws.Range("E45").Value = '  +4.30%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.19'
$ws.Range("E46").Value = '  -2.67%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '25.10'
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '7.19'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '23.30'
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '0.941'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '0.235'
$ws.Range("E51").Value = '  -1.48%  '

# Drop the temporary text number-format so the style index matches the
# original (unstyled) cells again.
$ws.Range("D2:D51").ClearFormats()
